$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.01244
$ws.Range("D2").Value = 0.8104
$ws.Range("E2").Value = 0.4966
$ws.Range("F2").Value = 11
$ws.Range("C3").Value = 0.1009
$ws.Range("D3").Value = 0.9441000000000001
$ws.Range("E3").Value = 0.09044000000000001
$ws.Range("F3").Value = 17
$ws.Range("C4").Value = 0.3462
$ws.Range("D4").Value = 1.275
$ws.Range("E4").Value = 0.2749
$ws.Range("F4").Value = 19
$ws.Range("C5").Value = 0.3308
$ws.Range("D5").Value = 0.968
$ws.Range("E5").Value = 0.1534
$ws.Range("F5").Value = 24
$ws.Range("C6").Value = 0.1351
$ws.Range("D6").Value = 1.069
$ws.Range("E6").Value = 0.1706
$ws.Range("F6").Value = 24
$ws.Range("C7").Value = 0.1689
$ws.Range("D7").Value = 1.459
$ws.Range("E7").Value = 0.4461
$ws.Range("F7").Value = 21
$ws.Range("C8").Value = 0.05829
$ws.Range("D8").Value = 1.546
$ws.Range("E8").Value = 1.399
$ws.Range("F8").Value = 14
$ws.Range("C9").Value = 0.1939
$ws.Range("D9").Value = 2.099
$ws.Range("E9").Value = 0.2264
$ws.Range("F9").Value = 16
$ws.Range("C10").Value = 0.4062
$ws.Range("D10").Value = 1.925
$ws.Range("E10").Value = 0.2214
$ws.Range("F10").Value = 22
$ws.Range("C11").Value = 0.4228
$ws.Range("D11").Value = 1.793
$ws.Range("E11").Value = 0.1071
$ws.Range("F11").Value = 24
$ws.Range("C12").Value = 0.0965
$ws.Range("D12").Value = 1.893
$ws.Range("E12").Value = 0.2045
$ws.Range("F12").Value = 11
$ws.Range("C13").Value = 0.002187
$ws.Range("D13").Value = 2.011
$ws.Range("E13").Value = 0.1693
$ws.Range("F13").Value = 9
$ws.Range("C14").Value = 0.1009
$ws.Range("D14").Value = 2.41
$ws.Range("E14").Value = 0.3006
$ws.Range("F14").Value = 14
$ws.Range("C15").Value = 0.09429999999999999
$ws.Range("D15").Value = 2.553
$ws.Range("E15").Value = 0.242
$ws.Range("F15").Value = 14
$ws.Range("C16").Value = 0.1088
$ws.Range("D16").Value = 2.786
$ws.Range("E16").Value = 0.5896
$ws.Range("F16").Value = 17
$ws.Range("C17").Value = 0.0002662
$ws.Range("D17").Value = 2.895
$ws.Range("E17").Value = 0.2886
$ws.Range("F17").Value = 9
$ws.Range("C18").Value = 0.007861
$ws.Range("D18").Value = 3.097
$ws.Range("E18").Value = 0.2737
$ws.Range("F18").Value = 14
$ws.Range("C19").Value = 0.1738
$ws.Range("D19").Value = 2.905
$ws.Range("E19").Value = 0.2202
$ws.Range("F19").Value = 24
$ws.Range("C20").Value = 0.01252
$ws.Range("D20").Value = 3.328
$ws.Range("E20").Value = 0.1245
$ws.Range("F20").Value = 17
$ws.Range("C21").Value = 0.0126
$ws.Range("D21").Value = 3.479
$ws.Range("E21").Value = 0.107
$ws.Range("F21").Value = 19
$ws.Range("C22").Value = 0.001719
$ws.Range("D22").Value = 3.94
$ws.Range("E22").Value = 0.105
$ws.Range("F22").Value = 12
$ws.Range("C23").Value = 0.0002045
$ws.Range("D23").Value = 4.967
$ws.Range("E23").Value = 0.3142
$ws.Range("F23").Value = 7
$ws.Range("C24").Value = 0.004486
$ws.Range("D24").Value = 4.275
$ws.Range("E24").Value = 0.4583
$ws.Range("F24").Value = 16
$ws.Range("C25").Value = 0.001384
$ws.Range("D25").Value = 5.265
$ws.Range("E25").Value = 0.3796
$ws.Range("F25").Value = 7
$ws.Range("C26").Value = 0.0008876999999999999
$ws.Range("D26").Value = 5.509
$ws.Range("E26").Value = 0.1375
$ws.Range("F26").Value = 7
$ws.Range("C27").Value = 0.006732
$ws.Range("D27").Value = 4.791
$ws.Range("E27").Value = 0.2124
$ws.Range("F27").Value = 12
$ws.Range("C28").Value = 0.01462
$ws.Range("D28").Value = 5.022
$ws.Range("E28").Value = 0.06677
$ws.Range("F28").Value = 16
$ws.Range("C29").Value = 0.0009169
$ws.Range("D29").Value = 4.787
$ws.Range("E29").Value = 0.04599
$ws.Range("F29").Value = 13
$ws.Range("C30").Value = 0.0001061
$ws.Range("D30").Value = 5.315
$ws.Range("E30").Value = 0.07238
$ws.Range("F30").Value = 7
$ws.Range("C31").Value = 0.06807000000000001
$ws.Range("D31").Value = 4.74
$ws.Range("E31").Value = 0.253
$ws.Range("F31").Value = 24
$ws.Range("C32").Value = 0.000009308
$ws.Range("D32").Value = 5.487
$ws.Range("E32").Value = 0.3387
$ws.Range("F32").Value = 7
$ws.Range("C33").Value = 0.01389
$ws.Range("D33").Value = 5.046
$ws.Range("E33").Value = 0.2695
$ws.Range("F33").Value = 14
$ws.Range("C34").Value = 0.001119
$ws.Range("D34").Value = 4.833
$ws.Range("E34").Value = 0.04521
$ws.Range("F34").Value = 12
$ws.Range("C35").Value = 0.001326
$ws.Range("D35").Value = 4.945
$ws.Range("E35").Value = 0.164
$ws.Range("F35").Value = 14
$ws.Range("C36").Value = 0.00000000000000000000000000002998
$ws.Range("D36").Value = 6.24
$ws.Range("E36").Value = 0.1078
$ws.Range("F36").Value = 6
$ws.Range("C37").Value = 0.05605
$ws.Range("D37").Value = 5.25
$ws.Range("E37").Value = 0.0853
$ws.Range("F37").Value = 22
$ws.Range("C38").Value = 0.001069
$ws.Range("D38").Value = 5.933
$ws.Range("E38").Value = 0.1564
$ws.Range("F38").Value = 7
$ws.Range("C39").Value = 0.004094
$ws.Range("D39").Value = 5.513
$ws.Range("E39").Value = 0.5288
$ws.Range("F39").Value = 12
$ws.Range("C40").Value = 0.01848
$ws.Range("D40").Value = 5.432
$ws.Range("E40").Value = 0.2089
$ws.Range("F40").Value = 14
$ws.Range("C41").Value = 0.006495
$ws.Range("D41").Value = 5.426
$ws.Range("E41").Value = 0.08316999999999999
$ws.Range("F41").Value = 14
$ws.Range("C42").Value = 0.006423
$ws.Range("D42").Value = 5.554
$ws.Range("E42").Value = 0.2838
$ws.Range("F42").Value = 16
$ws.Range("C43").Value = 0.1434
$ws.Range("D43").Value = 5.236
$ws.Range("E43").Value = 0.2819
$ws.Range("F43").Value = 24
$ws.Range("C44").Value = 0.0000000000000000000000000001231
$ws.Range("D44").Value = 6.423
$ws.Range("E44").Value = 0.7166
$ws.Range("F44").Value = 6
$ws.Range("C45").Value = 0.048
$ws.Range("D45").Value = 5.488
$ws.Range("E45").Value = 1.173
$ws.Range("F45").Value = 11
$ws.Range("C46").Value = 0.1221
$ws.Range("D46").Value = 5.687
$ws.Range("E46").Value = 0.4614
$ws.Range("F46").Value = 14
$ws.Range("C47").Value = 0.3026
$ws.Range("D47").Value = 5.804
$ws.Range("E47").Value = 0.07955
$ws.Range("F47").Value = 21
$ws.Range("C48").Value = 0.02337
$ws.Range("D48").Value = 5.823
$ws.Range("E48").Value = 0.06987
$ws.Range("F48").Value = 11
$ws.Range("C49").Value = 0.02265
$ws.Range("D49").Value = 5.867
$ws.Range("E49").Value = 0.4875
$ws.Range("F49").Value = 12
$ws.Range("C50").Value = 0.03538
$ws.Range("D50").Value = 5.914
$ws.Range("E50").Value = 1.27
$ws.Range("F50").Value = 12
$ws.Range("C51").Value = 0.003351
$ws.Range("D51").Value = 6.443
$ws.Range("E51").Value = 0.6151
$ws.Range("F51").Value = 7
$ws.Range("C52").Value = 0.03072
$ws.Range("D52").Value = 5.697
$ws.Range("E52").Value = 0.4801
$ws.Range("F52").Value = 11
$ws.Range("C53").Value = 0.08151
$ws.Range("D53").Value = 5.87
$ws.Range("E53").Value = 0.3352
$ws.Range("F53").Value = 19
$ws.Range("C54").Value = 0.6843
$ws.Range("D54").Value = 5.576
$ws.Range("E54").Value = 0.4235
$ws.Range("F54").Value = 24
$ws.Range("C55").Value = 0.008108000000000001
$ws.Range("D55").Value = 4.895
$ws.Range("E55").Value = 0.1901
$ws.Range("F55").Value = 11
$ws.Range("C56").Value = 0.07369000000000001
$ws.Range("D56").Value = 5.145
$ws.Range("E56").Value = 0.1319
$ws.Range("F56").Value = 14
$ws.Range("C57").Value = 0.007827000000000001
$ws.Range("D57").Value = 4.831
$ws.Range("E57").Value = 0.5063
$ws.Range("F57").Value = 11
$ws.Range("C58").Value = 0.01284
$ws.Range("D58").Value = 4.479
$ws.Range("E58").Value = 0.361
$ws.Range("F58").Value = 9
$ws.Range("C59").Value = 0.2662
$ws.Range("D59").Value = 5.16
$ws.Range("E59").Value = 0.09773999999999999
$ws.Range("F59").Value = 22
$ws.Range("C60").Value = 0.02578
$ws.Range("D60").Value = 4.62
$ws.Range("E60").Value = 0.07553
$ws.Range("F60").Value = 14
$ws.Range("C61").Value = 0.04162
$ws.Range("D61").Value = 4.699
$ws.Range("E61").Value = 0.124
$ws.Range("F61").Value = 17
$ws.Range("C62").Value = 0.1212
$ws.Range("D62").Value = 4.882
$ws.Range("E62").Value = 0.1589
$ws.Range("F62").Value = 22
$ws.Range("C63").Value = 0.1728
$ws.Range("D63").Value = 4.809
$ws.Range("E63").Value = 0.9268999999999999
$ws.Range("F63").Value = 24
$ws.Range("C64").Value = 0.006985
$ws.Range("D64").Value = 4.325
$ws.Range("E64").Value = 0.2466
$ws.Range("F64").Value = 11
$ws.Range("C65").Value = 0.004142
$ws.Range("D65").Value = 5.26
$ws.Range("E65").Value = 0.2186
$ws.Range("F65").Value = 7
$ws.Range("C66").Value = 0.01208
$ws.Range("D66").Value = 5.06
$ws.Range("E66").Value = 1.013
$ws.Range("F66").Value = 16
$ws.Range("C67").Value = 0.01316
$ws.Range("D67").Value = 5.189
$ws.Range("E67").Value = 0.177
$ws.Range("F67").Value = 16
$ws.Range("C68").Value = 0.08619
$ws.Range("D68").Value = 5.176
$ws.Range("E68").Value = 0.1181
$ws.Range("F68").Value = 21
$ws.Range("C69").Value = 0.006925
$ws.Range("D69").Value = 5.872
$ws.Range("E69").Value = 0.1013
$ws.Range("F69").Value = 9
$ws.Range("C70").Value = 0.02144
$ws.Range("D70").Value = 5.9
$ws.Range("E70").Value = 0.0873
$ws.Range("F70").Value = 14
$ws.Range("C71").Value = 0.01113
$ws.Range("D71").Value = 6.014
$ws.Range("E71").Value = 0.2145
$ws.Range("F71").Value = 12
$ws.Range("C72").Value = 0.02255
$ws.Range("D72").Value = 6.167
$ws.Range("E72").Value = 0.1816
$ws.Range("F72").Value = 17
